$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 74.72141563250115
    "C2" = 97.13589999711806
    "D2" = 99.38462181770601
    "E2" = 98.88086706005839
    "F2" = 98.24239346626013
    "G2" = 97.31529489657784
    "H2" = 95.84200289885352

    "B3" = 42.79607051046356
    "C3" = 97.55769539225048
    "D3" = 99.78788689456675
    "E3" = 99.04480641882745
    "F3" = 98.42782290177972
    "G3" = 97.46403351055639
    "H3" = 95.91340223752626

    "B4" = 91.1959470467071
    "C4" = 96.78665432045534
    "D4" = 99.22175657692412
    "E4" = 98.70586128554535
    "F4" = 98.32484031176851
    "G4" = 97.33507131787675
    "H4" = 95.82661924094829

    "B5" = 66.43068031500303
    "C5" = 96.87405924577853
    "D5" = 99.43339868583621
    "E5" = 98.81460183363905
    "F5" = 98.28709609094793
    "G5" = 97.16781939710914
    "H5" = 95.96845013475753

    "B6" = 69.49300814699041
    "C6" = 97.28585912329721
    "D6" = 99.35999213418873
    "E6" = 98.87461661348453
    "F6" = 98.2307537224838
    "G6" = 97.38278819347511
    "H6" = 95.91903370780105
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
